$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.311.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.930.33"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7361"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3214"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.82"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07080"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7856"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08029"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.932.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.381"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.55"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.59"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.300.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "254.33"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008040"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.739"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.184.74"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.822"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.557"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.80"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.298"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.46%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1321"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.359"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.535"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.427"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.147"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05118"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.287"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7467"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.773"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.804"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.04"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.401"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4502"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.984"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8433"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.39"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.534"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "974.23"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.87"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4183"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.92%  "
